$wb = $excel.ActiveWorkbook

# --- Update shared-string text values ---
$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# URL (StructureDefinition) appears both on the Metadata sheet (B2) and
# as the "Fixed Value" of the Elements sheet (Q5) - keep both in sync.
$meta.Range("B2").Value = "https://hl7chile.cl/fhir/ig/CoreCL/StructureDefinition/NombreCampana"
$elements.Range("Q5").Value = "https://hl7chile.cl/fhir/ig/CoreCL/StructureDefinition/NombreCampana"

# Date (Metadata sheet, B8)
$meta.Range("B8").Value = "2022-12-12T20:08:16-03:00"

# ValueSet URL ("Binding Value Set" on the Elements sheet, Y7)
$elements.Range("Y7").Value = "https://hl7chile.cl/fhir/ig/CoreCL/ValueSet/VSNombreCampana"

# --- Update (shrink) column widths on the Elements sheet, mirroring the
#     "best fit" recompute that happened when the workbook was re-saved ---
$elements.Columns.Item(1).ColumnWidth = 18.166666666666668
$elements.Columns.Item(2).ColumnWidth = 11.166666666666666
$elements.Columns.Item(3).ColumnWidth = 6.833333333333333
$elements.Columns.Item(4).ColumnWidth = 5.0
$elements.Columns.Item(5).ColumnWidth = 3.8333333333333335
$elements.Columns.Item(6).ColumnWidth = 4.166666666666667
$elements.Columns.Item(7).ColumnWidth = 13.833333333333334
$elements.Columns.Item(8).ColumnWidth = 11.166666666666666
$elements.Columns.Item(9).ColumnWidth = 11.833333333333334
$elements.Columns.Item(11).ColumnWidth = 89.16666666666667
$elements.Columns.Item(15).ColumnWidth = 12.666666666666666
$elements.Columns.Item(20).ColumnWidth = 14.833333333333334
$elements.Columns.Item(21).ColumnWidth = 15.333333333333334
$elements.Columns.Item(22).ColumnWidth = 16.166666666666668
$elements.Columns.Item(23).ColumnWidth = 15.5
$elements.Columns.Item(24).ColumnWidth = 18.0
$elements.Columns.Item(25).ColumnWidth = 57.166666666666664
$elements.Columns.Item(26).ColumnWidth = 4.833333333333333
$elements.Columns.Item(27).ColumnWidth = 18.833333333333332
$elements.Columns.Item(28).ColumnWidth = 39.166666666666664
$elements.Columns.Item(29).ColumnWidth = 14.166666666666666
$elements.Columns.Item(30).ColumnWidth = 11.5
$elements.Columns.Item(31).ColumnWidth = 16.833333333333332
$elements.Columns.Item(32).ColumnWidth = 8.666666666666666
$elements.Columns.Item(33).ColumnWidth = 9.0
$elements.Columns.Item(34).ColumnWidth = 11.333333333333334
$elements.Columns.Item(36).ColumnWidth = 21.833333333333332

# Columns 3, 4, 31, 32, 33 stayed hidden in the target workbook.
$elements.Columns.Item(3).Hidden = $true
$elements.Columns.Item(4).Hidden = $true
$elements.Columns.Item(31).Hidden = $true
$elements.Columns.Item(32).Hidden = $true
$elements.Columns.Item(33).Hidden = $true

$wb.Save()
